$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column C (13-01-2023), matching the style of B1
$ws.Range("C1").Value = "13-01-2023"
$ws.Range("B1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# Rebuild rows 2-10 with new labels/order and values for columns A, B, C
$data = @(
    @("Alpha Acciones",        38825.95, 39026.29),
    @("Alpha Mega",             5749.75,  5743.12),
    @("Arpenta acciones",       3240.52,  3246.14),
    @("Compass Crecimiento",          0,        0),
    @("Fima Acciones",                0,        0),
    @("Fima PB Acciones",             0,        0),
    @("HF Acciones Argentinas", 1121.59,  1109.86),
    @("avg",                    6991.12,  7017.92),
    @("total",                 48937.81, 49125.41)
)

$row = 2
foreach ($item in $data) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = $item[2]
    $row++
}
